$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Random Forest
$ws.Range("G2").Value = 99.62
$ws.Range("H2").Value = 0.77

# Row 3 - SVM
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 100

# Row 4 - KNN
$ws.Range("G4").Value = 98.84999999999999
$ws.Range("H4").Value = 1.54

# Row 5 - Regressão Logística
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 100
$ws.Range("H5").Value = 1.22
